$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header texts for columns A..AT (row 1), reflecting the rebuilt shared string table.
$headers = @(
    "CODIGO DE PROYECTO",
    "CODIGO DE PRODUCTO",
    "NOMBRE DEL PRODUCTO",
    "DIRECTORIO",
    "AÑO",
    "AUTOR 1",
    "AUTOR 2",
    "AUTOR 3",
    "AUTOR 4",
    "AUTOR 5",
    "DESCRIPCION",
    "TIPO DE PRODUCTO",
    "TIPO DE MODELO",
    "AREA DE CONOCIMIENTO DEL MODELO",
    "SE CUENTA CON ARCHIVOS NATIVOS",
    "SOFTWARE",
    "TIPO DE LICENCIA",
    "NIVEL DE CONSUMO DE RECURSOS COMPUTACIONALES",
    "CODIGO PRODUCTO ASOCIADO 1",
    "CODIGO PRODUCTO ASOCIADO 2",
    "CODIGO PRODUCTO ASOCIADO 3",
    "CODIGO PRODUCTO ASOCIADO 4",
    "CODIGO PRODUCTO ASOCIADO 5",
    "CANTIDAD DE ESCENARIOS SIMULADOS",
    "NOMBRE ESCENARIO 1",
    "DESCRIPCION ESCENARIO 1",
    "NOMBRE ESCENARIO 2",
    "DESCRIPCION ESCENARIO 2",
    "NOMBRE ESCENARIO 3",
    "DESCRIPCION ESCENARIO 3",
    "NOMBRE ESCENARIO 4",
    "DESCRIPCION ESCENARIO 4",
    "NOMBRE ESCENARIO 5",
    "DESCRIPCION ESCENARIO 5",
    "NOMBRE ESCENARIO 6",
    "DESCRIPCION ESCENARIO 6",
    "NOMBRE ESCENARIO 7",
    "DESCRIPCION ESCENARIO 7",
    "NOMBRE ESCENARIO 8",
    "DESCRIPCION ESCENARIO 8",
    "NOMBRE ESCENARIO 9",
    "DESCRIPCION ESCENARIO 9",
    "NOMBRE ESCENARIO 10",
    "DESCRIPCION ESCENARIO 10",
    "OBSERVACIONES",
    "PALABRAS CLAVE"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = $i + 1
    $cell = $ws.Cells.Item(1, $col)
    $cell.Value = $headers[$i]
}

# Extend the bordered, empty style block for rows 1-4 out to column AT (46),
# matching the same formatting already used in columns A-AH (a thin black
# border around every cell).
$newRange = $ws.Range("AI1:AT4")
$newRange.Borders.ColorIndex = 1
$newRange.Borders.LineStyle = 1
$newRange.Borders.Weight = 2

# Final selection matches the saved view: cell J12 selected (also drops the
# old topLeftCell="W1" scroll anchor, same as the target workbook).
$ws.Range("J12").Select() | Out-Null
